$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.287.72'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '2.695.05'
$ws.Range("E3").Value = '  +1.61%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.97'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.24'
$ws.Range("E6").Value = '  +1.98%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.127'
$ws.Range("E9").Value = '  +8.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.02'
$ws.Range("E10").Value = '  +3.47%  '
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000211'
$ws.Range("E13").Value = '  +21.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.36'
$ws.Range("E14").Value = '  +3.80%  '
$ws.Range("D15").Value = '3.177.38'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").Value = '66.103.65'
$ws.Range("E16").Value = '  +1.31%  '
$ws.Range("D17").Value = '2.689.42'
$ws.Range("E17").Value = '  +1.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.74'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.70'
$ws.Range("E20").Value = '  +2.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.57'
$ws.Range("E21").Value = '  +4.04%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.15'
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.77'
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000108'
$ws.Range("E25").Value = '  +16.52%  '
$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("E26").Value = '  -2.67%  '
$ws.Range("E27").Value = '  +5.79%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.24'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("E30").Value = '  +7.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '541.14'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("E34").Value = '  +2.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").Value = '  -5.77%  '
$ws.Range("E36").Value = '  +1.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.87'
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.48'
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.02'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '171.18'
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.75'
$ws.Range("E43").Value = '  +2.46%  '
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.36'
$ws.Range("E45").Value = '  +3.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0622'
$ws.Range("E46").Value = '  +2.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.35'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.666'
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("E49").Value = '  +6.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.48'
$ws.Range("E50").Value = '  +4.69%  '
$ws.Range("E51").Value = '  +0.43%  '
